$d = $word.ActiveDocument

# --- Insert a new underlined "TP0" paragraph at the very start of the document ---
$firstRange = $d.Paragraphs.Item(1).Range
$firstRange.InsertParagraphBefore()
$tp0 = $d.Paragraphs.Item(1)
$tp0.Range.Text = "TP0"
$tp0.Range.Font.Underline = 1

# --- Locate the "No additional modules planned to be used" paragraph ---
$idx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*No additional modules planned to be used*") {
        $idx = $i
    }
}

# --- Append two new paragraphs after it: "TP1 Update" (underlined) and the update text ---
$target = $d.Paragraphs.Item($idx)
$target.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($idx + 1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($idx + 2)

$p1.Range.Text = "TP1 Update"
$p2.Range.Text = "Project will only consist of one player mode and competitive mode. Project will also have rats spawn throughout the game in a random location. These rats will target the food that was last left on the counter, if there is food on the counter, it will move towards the food and steal it unless the player moves the food away."

$p1.Range.Font.Underline = 1
